$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new values for columns C (nombre_aides) and E (montant_total)
$updates = @{
    3   = @{ C = 249328;  E = 1036478313 }
    53  = @{ C = 141680;  E = 590063206 }
    57  = @{ C = 3712;    E = 138408892 }
    92  = @{ C = 409074;  E = 1594755556 }
    93  = @{ C = 209550;  E = 1308407897 }
    95  = @{ C = 50755;   E = 931923076 }
    96  = @{ C = 17257;   E = 790651633 }
    104 = @{ C = 135233;  E = 272168047 }
    110 = @{ C = 396;     E = 16649846 }
    174 = @{ C = 226089;  E = 900577214 }
    175 = @{ C = 80780;   E = 486154029 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}
